$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New daily records (date serial, nuovi pos., somma mobile 7gg., somma mobile 7gg. per 100mila abitanti)
$data = @(
    @(375, 44449, 0, 3, 19.88466892026248),
    @(376, 44450, 0, 3, 19.88466892026248),
    @(377, 44451, 0, 3, 19.88466892026248),
    @(378, 44452, 0, 1, 6.628222973420826),
    @(379, 44453, 2, 3, 19.88466892026248),
    @(380, 44454, 1, 4, 26.5128918936833),
    @(381, 44455, 2, 5, 33.14111486710413),
    @(382, 44456, 0, 5, 33.14111486710413),
    @(383, 44457, 1, 6, 39.76933784052495),
    @(384, 44458, 1, 7, 46.39756081394578),
    @(385, 44459, 3, 10, 66.28222973420826)
)

foreach ($row in $data) {
    $r = $row[0]
    $ws.Cells.Item($r, 1).Value = $row[1]
    $ws.Cells.Item($r, 2).Value = $row[2]
    $ws.Cells.Item($r, 3).Value = $row[3]
    $ws.Cells.Item($r, 4).Value = $row[4]
}

# Match the formatting (date number format, font, border, alignment) used by the
# existing data rows, carrying it down onto each newly appended row.
$ws.Range("A374:D374").Copy()
foreach ($row in $data) {
    $r = $row[0]
    $ws.Range("A" + $r + ":D" + $r).PasteSpecial(-4122)
}

Write-Output "Updated through row 385 (2021-09-20)"
